# Update the "Input" worksheet's dummy template values:
#  - replace the sample server URL with a placeholder URL
#  - replace the per-row project names with generic placeholders
#  - replace the real-looking PAT with a dummy PAT
#  - drop the last two sample rows (only two example rows remain)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Remember the current (non-hyperlink) look of column A so we can restore it
# after re-creating the hyperlinks below (Hyperlinks.Add applies the builtin
# "Hyperlink" style, but these cells already use a custom look-alike style).
$origStyleA2 = $ws.Range("A2").Style
$origStyleA3 = $ws.Range("A3").Style

# Drop all the existing hyperlinks up front; we'll recreate just the two we
# keep once the final row layout is in place.
$ws.Hyperlinks.Delete()

# Row 2: project1
$ws.Range("A2").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B2").Value = "project1"
$ws.Range("C2").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# Row 3: project2
$ws.Range("A3").Value = "http://128.0.0.1/TestCollection"
$ws.Range("B3").Value = "project2"
$ws.Range("C3").Value = "adad87adad8ds4449m434344mmnbnbb43434"

# Remove the old rows 4 and 5 (qaserver/testproj samples), shifting cells up.
$ws.Range("A4:C5").Delete()

# Recreate the hyperlinks for the remaining two rows against the new URL.
$ws.Hyperlinks.Add($ws.Range("A2"), "http://128.0.0.1/TestCollection")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://128.0.0.1/TestCollection")

# Restore the original cell look now that the hyperlink has been (re)applied.
$ws.Range("A2").Style = $origStyleA2
$ws.Range("A3").Style = $origStyleA3
